# Auto-generated script applying the cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.875.19"
$ws.Range("E2").Value = "  -0.55%  "
$ws.Range("D3").Value = "'2.521.88"
$ws.Range("E3").Value = "  +0.33%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").Value = "'536.02"
$ws.Range("E5").Value = "  +0.31%  "
$ws.Range("D6").Value = "'137.54"
$ws.Range("E6").Value = "  -1.13%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "'0.566"
$ws.Range("E8").Value = "  +0.30%  "
$ws.Range("D9").Value = "'2.520.80"
$ws.Range("E9").Value = "  +0.27%  "
$ws.Range("E10").Value = "  +0.42%  "
$ws.Range("E11").Value = "  -2.08%  "
$ws.Range("E12").Value = "  -1.12%  "
$ws.Range("D13").Value = "'0.348"
$ws.Range("E13").Value = "  -2.15%  "
$ws.Range("D14").Value = "'2.969.65"
$ws.Range("E14").Value = "  +0.23%  "
$ws.Range("D15").Value = "'23.02"
$ws.Range("E15").Value = "  -1.50%  "
$ws.Range("D16").Value = "'58.895.32"
$ws.Range("E16").Value = "  -0.39%  "
$ws.Range("E17").Value = "  -0.96%  "
$ws.Range("D18").Value = "'2.520.98"
$ws.Range("E18").Value = "  +0.15%  "
$ws.Range("D19").Value = "'11.11"
$ws.Range("E19").Value = "  +0.21%  "
$ws.Range("E20").Value = "  -0.31%  "
$ws.Range("D21").Value = "'325.55"
$ws.Range("E21").Value = "  +0.42%  "
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").Value = "'5.97"
$ws.Range("E23").Value = "  +2.95%  "
$ws.Range("D24").Value = "'65.59"
$ws.Range("E24").Value = "  +2.85%  "
$ws.Range("D25").Value = "'0.422"
$ws.Range("E25").Value = "  -1.17%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("D28").Value = "'7.63"
$ws.Range("E28").Value = "  -2.39%  "
$ws.Range("D29").Value = "'6.71"
$ws.Range("E29").Value = "  -2.48%  "
$ws.Range("D30").Value = "'0.0₃0768"
$ws.Range("E30").Value = "  -0.70%  "
$ws.Range("D31").Value = "'1.78"
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("E32").Value = "  +6.41%  "
$ws.Range("D33").Value = "'163.79"
$ws.Range("E33").Value = "  -0.26%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "'1.47"
$ws.Range("E34").Value = "  +0.08%  "
$ws.Range("B35").Value = "USDe"
$ws.Range("C35").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D35").Value = "'0.998"
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("E36").Value = "  -0.16%  "
$ws.Range("E37").Value = "  -3.04%  "
$ws.Range("D39").Value = "'36.67"
$ws.Range("E39").Value = "  -0.63%  "
$ws.Range("D40").Value = "'0.819"
$ws.Range("E40").Value = "  +0.65%  "
$ws.Range("E41").Value = "  -1.28%  "
$ws.Range("D42").Value = "'285.98"
$ws.Range("E42").Value = "  +2.72%  "
$ws.Range("D43").Value = "'5.23"
$ws.Range("E43").Value = "  -0.05%  "
$ws.Range("D44").Value = "'132.68"
$ws.Range("E44").Value = "  +7.48%  "
$ws.Range("D45").Value = "'0.997"
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("D46").Value = "'0.604"
$ws.Range("E46").Value = "  +1.40%  "
$ws.Range("D47").Value = "'10.88"
$ws.Range("E47").Value = "  +0.14%  "
$ws.Range("D48").Value = "'0.0930"
$ws.Range("E48").Value = "  -0.43%  "
$ws.Range("D49").Value = "'0.0509"
$ws.Range("E49").Value = "  -0.53%  "
$ws.Range("E50").Value = "  -1.14%  "
$ws.Range("D51").Value = "'17.29"
$ws.Range("E51").Value = "  -2.26%  "
